# Update EUR->ARS rate: append new reading as row 33
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the date/time-looking strings to be stored as literal text
# (matching the existing rows, which are plain inline/shared strings,
# not real Excel dates), then restore the default "Normal" style so no
# stray number-format styling is left behind on the new cells.
$ws.Range("A33").NumberFormat = "@"
$ws.Range("A33").Value = "2025-09-22"
$ws.Range("A33").Style = "Normal"

$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "15:23:16"
$ws.Range("B33").Style = "Normal"

$ws.Range("C33").Value = "1.00 EUR = 1,782.7350"
